$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = "59.353.93"
$cD.Style = "Normal"
$ws.Range("E2").Value = "  +2.70%  "

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = "2.981.94"
$cD.Style = "Normal"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  -0.02%  "

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = "561.86"
$cD.Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "

$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value = "138.16"
$cD.Style = "Normal"
$ws.Range("E6").Value = "  +4.88%  "

$ws.Range("E7").Value = "  -0.10%  "

$cD = $ws.Range("D8")
$cD.NumberFormat = "@"
$cD.Value = "0.518"
$cD.Style = "Normal"

$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = "2.972.49"
$cD.Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "

$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value = "0.132"
$cD.Style = "Normal"
$ws.Range("E10").Value = "  +4.08%  "

$ws.Range("E11").Value = "  +11.54%  "

$ws.Range("E12").Value = "  +1.81%  "

$cD = $ws.Range("D13")
$cD.NumberFormat = "@"
$cD.Value = "0.0000229"
$cD.Style = "Normal"
$ws.Range("E13").Value = "  +4.11%  "

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = "33.67"
$cD.Style = "Normal"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("E15").Value = "  -0.54%  "

$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value = "3.482.21"
$cD.Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "

$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = "7.14"
$cD.Style = "Normal"
$ws.Range("E17").Value = "  +4.06%  "

$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = "2.989.47"
$cD.Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = "59.375.48"
$cD.Style = "Normal"
$ws.Range("E19").Value = "  +2.73%  "

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = "433.32"
$cD.Style = "Normal"
$ws.Range("E20").Value = "  +4.29%  "

$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = "13.55"
$cD.Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "

$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = "0.718"
$cD.Style = "Normal"
$ws.Range("E22").Value = "  +3.48%  "

$cD = $ws.Range("D23")
$cD.NumberFormat = "@"
$cD.Value = "13.33"
$cD.Style = "Normal"
$ws.Range("E23").Value = "  -1.13%  "

$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value = "7.00"
$cD.Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value = "79.76"
$cD.Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +10.37%  "

$ws.Range("E28").Value = "  -0.07%  "

$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value = "7.71"
$cD.Style = "Normal"
$ws.Range("E30").Value = "  +4.67%  "

$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = "0.106"
$cD.Style = "Normal"
$ws.Range("E31").Value = "  +8.56%  "

$cD = $ws.Range("D32")
$cD.NumberFormat = "@"
$cD.Value = "6.24"
$cD.Style = "Normal"
$ws.Range("E32").Value = "  +4.97%  "

$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = "25.68"
$cD.Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "

$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = "0.0₃0773"
$cD.Style = "Normal"
$ws.Range("E34").Value = "  +11.25%  "

$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value = "0.990"
$cD.Style = "Normal"
$ws.Range("E35").Value = "  +5.69%  "

$cD = $ws.Range("D36")
$cD.NumberFormat = "@"
$cD.Value = "5.86"
$cD.Style = "Normal"
$ws.Range("E36").Value = "  +3.81%  "

$cD = $ws.Range("D37")
$cD.NumberFormat = "@"
$cD.Value = "2.07"
$cD.Style = "Normal"
$ws.Range("E37").Value = "  +0.41%  "

$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = "48.75"
$cD.Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "

$cD = $ws.Range("D39")
$cD.NumberFormat = "@"
$cD.Value = "8.65"
$cD.Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "

$cD = $ws.Range("D40")
$cD.NumberFormat = "@"
$cD.Value = "2.77"
$cD.Style = "Normal"
$ws.Range("E40").Value = "  +6.42%  "

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = "398.86"
$cD.Style = "Normal"
$ws.Range("E41").Value = "  +6.45%  "

$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value = "0.0352"
$cD.Style = "Normal"
$ws.Range("E42").Value = "  +2.25%  "

$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = "2.750.70"
$cD.Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "

$ws.Range("E44").Value = "  -2.44%  "

$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = "0.250"
$cD.Style = "Normal"
$ws.Range("E45").Value = "  +6.09%  "

$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value = "122.71"
$cD.Style = "Normal"
$ws.Range("E47").Value = "  -1.00%  "

$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value = "34.34"
$cD.Style = "Normal"
$ws.Range("E48").Value = "  +18.15%  "

$ws.Range("E49").Value = "  +1.80%  "

$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value = "2.00"
$cD.Style = "Normal"
$ws.Range("E50").Value = "  +2.38%  "

$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = "23.31"
$cD.Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
